$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 105; $r -le 143; $r++) {
    $ws.Cells.Item($r, 1).Value = 1.147256941293532
}
